# Auto-generated edit script applying the diff changes to before.xlsx
# Updates three sheets: VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL
$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item('VENTAS POR GRUPO')
$wsGrupo.Range('D70').Value = 5280.57
$wsGrupo.Range('L70').Value = 506.88
$wsGrupo.Range('M70').Value = 2473.79
$wsGrupo.Range('B172').Value = 'ARMIJOS SALINAS LUIS CLAUDIO'
$wsGrupo.Range('B173').Value = 'ASES GAVILANEZ FAUSTO HERNAN'
$wsGrupo.Range('B174').Value = 'BARROS YUNGA DIEGO VINICIO'
$wsGrupo.Range('B175').Value = 'BRAVO MONTENEGRO DANIEL ANDRES'
$wsGrupo.Range('B176').Value = 'BRITO CARDENAS RUTH CECILIA'
$wsGrupo.Range('B177').Value = 'COELLO TRONCOSO JOSE GREGORIO'
$wsGrupo.Range('B178').Value = 'COMERCIAL LUNA PAZMIÑO CIA. LTDA.'
$wsGrupo.Range('M178').Value = 0
$wsGrupo.Range('B179').Value = 'CORPORACION AREVALO-YUMBLA E HIJOS'
$wsGrupo.Range('M179').Value = 1631.15
$wsGrupo.Range('B180').Value = 'FABIMP BENIGNO BRAVO S.A.S.'
$wsGrupo.Range('B181').Value = 'FRANK FERRETERIA FRANKFERRE CIA.'
$wsGrupo.Range('B182').Value = 'HUERTA MUÑOZ NANCY ELIZABETH'
$wsGrupo.Range('M182').Value = 0
$wsGrupo.Range('B183').Value = 'ILLER LOPEZ ROBERTO FERNANDO'
$wsGrupo.Range('M183').Value = 103.71
$wsGrupo.Range('B184').Value = 'MIM CONSTRUFERRETERIA E IMPORTADORA SAS'
$wsGrupo.Range('B185').Value = 'MOROCHO BACUILIMA HILDA INES'
$wsGrupo.Range('B186').Value = 'MULLO GUACHO ANA LUCIA'
$wsGrupo.Range('B187').Value = 'PAUTA ASTUDILLO JULIO HERNAN'
$wsGrupo.Range('L187').Value = 0
$wsGrupo.Range('B188').Value = 'ROCAFUERTE LOPEZ EVELYN ESTEFANIA'
$wsGrupo.Range('L188').Value = 591.61
$wsGrupo.Range('A189').Value = 'ILLER LOPEZ ROBERTO FERNANDO'
$wsGrupo.Range('B189').Value = 'VIEJO RIVAS MAYRA ANABELLE'
$wsGrupo.Range('M189').Value = 0
$wsGrupo.Range('B190').Value = 'ALARCON MAYORGA LUIS ALFONSO'
$wsGrupo.Range('M190').Value = 1097.67
$wsGrupo.Range('B191').Value = 'ALTAMIRANO ARIAS LUCIA ELIZABETH'
$wsGrupo.Range('B192').Value = 'ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO'
$wsGrupo.Range('B193').Value = 'ALVAREZ SAAVEDRA EDWIN GEOVANNY'
$wsGrupo.Range('B194').Value = 'ANGAMARCA CURIPONA WILMA'
$wsGrupo.Range('B195').Value = 'ARIAS MEZA RONALD FABRICIO'
$wsMensual = $wb.Worksheets.Item('VENTA MENSUAL')
$wsMensual.Range('F70').Value = 7840.56
$wsMensual.Range('B176').Value = 'ARMIJOS SALINAS LUIS CLAUDIO'
$wsMensual.Range('B177').Value = 'ASES GAVILANEZ FAUSTO HERNAN'
$wsMensual.Range('D177').Value = 0
$wsMensual.Range('E177').Value = 0
$wsMensual.Range('B178').Value = 'BARROS YUNGA DIEGO VINICIO'
$wsMensual.Range('D178').Value = 595.08
$wsMensual.Range('E178').Value = 17.99
$wsMensual.Range('B179').Value = 'BRAVO MONTENEGRO DANIEL ANDRES'
$wsMensual.Range('B180').Value = 'BRITO CARDENAS RUTH CECILIA'
$wsMensual.Range('B181').Value = 'COELLO TRONCOSO JOSE GREGORIO'
$wsMensual.Range('B182').Value = 'COMERCIAL LUNA PAZMIÑO CIA. LTDA.'
$wsMensual.Range('E182').Value = 0
$wsMensual.Range('F182').Value = 0
$wsMensual.Range('B183').Value = 'CORPORACION AREVALO-YUMBLA E HIJOS'
$wsMensual.Range('D183').Value = 0
$wsMensual.Range('E183').Value = 3992.9
$wsMensual.Range('F183').Value = 1631.15
$wsMensual.Range('B184').Value = 'FABIMP BENIGNO BRAVO S.A.S.'
$wsMensual.Range('D184').Value = 400.46
$wsMensual.Range('E184').Value = 252.25
$wsMensual.Range('B185').Value = 'FRANK FERRETERIA FRANKFERRE CIA.'
$wsMensual.Range('D185').Value = 5372.02
$wsMensual.Range('E185').Value = 7662.57
$wsMensual.Range('B186').Value = 'HUERTA MUÑOZ NANCY ELIZABETH'
$wsMensual.Range('D186').Value = 0
$wsMensual.Range('E186').Value = 0
$wsMensual.Range('F186').Value = 0
$wsMensual.Range('B187').Value = 'ILLER LOPEZ ROBERTO FERNANDO'
$wsMensual.Range('D187').Value = 136.48
$wsMensual.Range('E187').Value = 58.48
$wsMensual.Range('F187').Value = 103.71
$wsMensual.Range('B188').Value = 'MIM CONSTRUFERRETERIA E IMPORTADORA SAS'
$wsMensual.Range('D188').Value = 0
$wsMensual.Range('E188').Value = 3896.18
$wsMensual.Range('B189').Value = 'MOROCHO BACUILIMA HILDA INES'
$wsMensual.Range('D189').Value = 102.6
$wsMensual.Range('B190').Value = 'MULLO GUACHO ANA LUCIA'
$wsMensual.Range('E190').Value = 0
$wsMensual.Range('B191').Value = 'PAUTA ASTUDILLO JULIO HERNAN'
$wsMensual.Range('E191').Value = 326.73
$wsMensual.Range('F191').Value = 0
$wsMensual.Range('B192').Value = 'ROCAFUERTE LOPEZ EVELYN ESTEFANIA'
$wsMensual.Range('E192').Value = 1015.74
$wsMensual.Range('F192').Value = 591.61
$wsMensual.Range('A193').Value = 'ILLER LOPEZ ROBERTO FERNANDO'
$wsMensual.Range('B193').Value = 'VIEJO RIVAS MAYRA ANABELLE'
$wsMensual.Range('C193').Value = 0
$wsMensual.Range('D193').Value = 0
$wsMensual.Range('E193').Value = 0
$wsMensual.Range('F193').Value = 0
$wsMensual.Range('G193').Value = 0
$wsMensual.Range('B194').Value = 'ALARCON MAYORGA LUIS ALFONSO'
$wsMensual.Range('C194').Value = 1603.18
$wsMensual.Range('D194').Value = 2774.22
$wsMensual.Range('E194').Value = 1114.87
$wsMensual.Range('F194').Value = 1097.67
$wsMensual.Range('G194').Value = 2000
$wsMensual.Range('B195').Value = 'ALTAMIRANO ARIAS LUCIA ELIZABETH'
$wsMensual.Range('D195').Value = 1672.61
$wsMensual.Range('E195').Value = 6405.28
$wsMensual.Range('G195').Value = 2500
$wsMensual.Range('B196').Value = 'ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO'
$wsMensual.Range('E196').Value = 848.76
$wsMensual.Range('B197').Value = 'ALVAREZ SAAVEDRA EDWIN GEOVANNY'
$wsMensual.Range('B198').Value = 'ANGAMARCA CURIPONA WILMA'
$wsMensual.Range('C198').Value = 0
$wsMensual.Range('G198').Value = 0
$wsMensual.Range('B199').Value = 'ARIAS MEZA RONALD FABRICIO'
$wsMensual.Range('C199').Value = 2003.44
$wsMensual.Range('G199').Value = 1500
$wsMensual.Range('F362').Value = 291936.16
$wsCumplimiento = $wb.Worksheets.Item('CUMPLIMIENTO MENSUAL')
$wsCumplimiento.Range('D15').Value = 21657.75
$wsCumplimiento.Range('E15').Value = -6832.34
$wsCumplimiento.Range('F15').Value = 1.460853359198835
$wsCumplimiento.Range('D23').Value = 8165.68
$wsCumplimiento.Range('E23').Value = 7982.32
$wsCumplimiento.Range('F23').Value = 0.5056774832796631
$wsCumplimiento.Range('D24').Value = 22386.12
$wsCumplimiento.Range('E24').Value = 27920.88
$wsCumplimiento.Range('F24').Value = 0.4449901604150516
$wsCumplimiento.Range('D77').Value = 293829.05
$wsCumplimiento.Range('E77').Value = 123419.6297415455
$wsCumplimiento.Range('F77').Value = 0.7042060628735968
